# Implementing material and dynamic parameter rename an effect's texture
#
# A new row is inserted at row 122 of Sheet1 (the Languages resource sheet)
# holding the "DynamicParameter_Name" / "Dynamic parameter" / "動的パラメーター"
# localization triple, followed by a blank spacer row - mirroring the blank
# spacer rows that already separate the other localization groups further
# down the sheet. Everything that used to live at row 122 onward shifts
# down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 122 (shifts old row 122.. down to 124..),
# inheriting formatting/styles from the surrounding rows the same way
# Excel's own "Insert Row" does.
$ws.Rows.Item(122).Insert()
$ws.Rows.Item(122).Insert()

# Populate the new row 122 with the Dynamic Parameter strings.
# Set B and C before A so new shared-string entries are appended to
# sharedStrings.xml in the same order the source workbook uses them.
$ws.Range("B122").Value = "Dynamic parameter"
$ws.Range("C122").Value = "動的パラメーター"
$ws.Range("A122").Value = "DynamicParameter_Name"

# Row 123 is left as the blank separator row (already blank from the insert).

# Update the view: selection moves to A123, and the window scrolls so
# row 106 is the topmost visible row.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 106
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A123").Select()
